$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9206383228302002
$ws.Range("B1").Value = 1.732064127922058
$ws.Range("C1").Value = 4.115060329437256
$ws.Range("D1").Value = 3.59550952911377
$ws.Range("E1").Value = 0.3812944889068604
